# Activity Log workbook update:
#  - fill in rows 21-23 (log entries for 5-4-2020: final verification,
#    introduction writing, final report check)
#  - insert a blank formatted row above the old closing (thick-bottom) row
#    of the small table around row 54, pushing the table's bottom border
#    down one row
#  - move the active selection to where the user was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 21-23: new log entries (date column uses a leading apostrophe so
#     the text "5-4-2020" is stored as a string, matching the existing
#     date-formatted-as-text cells in column C rather than being parsed
#     into a real date serial) ---
$ws.Range("B21").Value = 779
$ws.Range("C21").Value = "'5-4-2020"
$ws.Range("D21").Value = 0.8354166666666667
$ws.Range("E21").Value = 0.85069444444444453
$ws.Range("G21").Value = "Final verification of our code"

$ws.Range("B22").Value = 779
$ws.Range("C22").Value = "'5-4-2020"
$ws.Range("D22").Value = 0.85069444444444453
$ws.Range("E22").Value = 0.92847222222222225

$ws.Range("B23").Value = 779
$ws.Range("C23").Value = "'5-4-2020"
$ws.Range("D23").Value = 0.92847222222222225
$ws.Range("E23").Value = 0.97361111111111109
$ws.Range("G23").Value = "Final check on report; redid Timing Logic Unit screenshots"

# (G22 written after G23 so new shared-string entries land in the same
#  order as the source: verification, final-check, introduction)
$ws.Range("G22").Value = "Wrote and revised introduction of report"

# --- Insert a row above the table's thick-bottom closing row (currently
#     row 54), shifting it (and everything below) down by one; then copy
#     the formatting of the row above into the newly inserted row so it
#     matches the rest of the table body ---
$ws.Rows("54:54").Insert()

$ws.Range("B53").Copy()
$ws.Range("B54").PasteSpecial(-4122)
$ws.Range("C53").Copy()
$ws.Range("C54").PasteSpecial(-4122)
$ws.Range("D53").Copy()
$ws.Range("D54").PasteSpecial(-4122)
$ws.Range("E53").Copy()
$ws.Range("E54").PasteSpecial(-4122)
$ws.Range("G53").Copy()
$ws.Range("G54").PasteSpecial(-4122)
$ws.Rows(54).RowHeight = $ws.Rows(53).RowHeight

$excel.CutCopyMode = $false

# --- Restore the cursor/selection to where the author left off ---
$ws.Range("D26").Select()
